$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.936.12"
$ws.Range("E2").Value = "  +4.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.804.29"
$ws.Range("E3").Value = "  +4.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -1.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "421.22"
$ws.Range("E5").Value = "  +4.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.58"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.798.53"
$ws.Range("E7").Value = "  +4.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.603"
$ws.Range("E8").Value = "  -2.38%  "
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.717"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.159"
$ws.Range("E11").Value = "  +1.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000347"
$ws.Range("E12").Value = "  +13.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.14"
$ws.Range("E13").Value = "  -3.58%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.415.00"
$ws.Range("E14").Value = "  +3.22%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10.12"
$ws.Range("E15").Value = "  +3.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.64"
$ws.Range("E16").Value = "  +22.10%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.137"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.788.15"
$ws.Range("E18").Value = "  +4.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.49"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.999.66"
$ws.Range("E20").Value = "  +2.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.07"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "405.26"
$ws.Range("E22").Value = "  -2.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.66"
$ws.Range("E23").Value = "  -2.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.52"
$ws.Range("E24").Value = "  -2.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.00"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.00"
$ws.Range("E26").Value = "  +4.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.50"
$ws.Range("E27").Value = "  +9.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.20"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.49"
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.04"
$ws.Range("E30").Value = "  +32.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "725.17"
$ws.Range("E31").Value = "  +7.88%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.43"
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.76"
$ws.Range("E33").Value = "  +2.67%  "
$ws.Range("E34").Value = "  +2.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.152"
$ws.Range("E36").Value = "  -4.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.56"
$ws.Range("E37").Value = "  -3.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.08"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.26"
$ws.Range("E39").Value = "  +24.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0739"
$ws.Range("E40").Value = "  +19.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0449"
$ws.Range("E41").Value = "  -2.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.91"
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("E44").Value = "  -3.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.32"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.315"
$ws.Range("E46").Value = "  +9.92%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.41"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.08"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.03"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.56"
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.78"
$ws.Range("E51").Value = "  +0.73%  "
